$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 93.62780766666667
$ws.Range("H2").Value = 280.883423
$ws.Range("I2").Value = 0.3228593149748609
$ws.Range("J2").Value = 0.3228593149748609
$ws.Range("M2").Value = 0.01046566666666667
$ws.Range("N2").Value = 0.031397
$ws.Range("O2").Value = 0.007547709568116775
$ws.Range("P2").Value = 0.007547709568116775
$ws.Range("Q2").Value = 0.9798774257701112
$ws.Range("R2").Value = 8.818896831930999
$ws.Range("S2").Value = 0.002436848340791385
$ws.Range("T2").Value = 0.002436848340791385
$ws.Range("G3").Value = 93.62780766666667
$ws.Range("H3").Value = 280.883423
$ws.Range("I3").Value = 0.3228593149748609
$ws.Range("J3").Value = 0.3228593149748609
$ws.Range("O3").Value = 0.8438809030711777
$ws.Range("P3").Value = 0.8438809030711776
$ws.Range("Q3").Value = 109.5563944923045
$ws.Range("R3").Value = 986.00755043074
$ws.Range("S3").Value = 0.2724548102859274
$ws.Range("T3").Value = 0.2724548102859273
$ws.Range("G4").Value = 93.62780766666667
$ws.Range("H4").Value = 280.883423
$ws.Range("I4").Value = 0.3228593149748609
$ws.Range("J4").Value = 0.3228593149748609
$ws.Range("M4").Value = 0.2060093333333333
$ws.Range("N4").Value = 0.618028
$ws.Range("O4").Value = 0.1485713873607056
$ws.Range("P4").Value = 0.1485713873607056
$ws.Range("Q4").Value = 19.28820223887156
$ws.Range("R4").Value = 173.593820149844
$ws.Range("S4").Value = 0.04796765634814211
$ws.Range("T4").Value = 0.04796765634814211
$ws.Range("G5").Value = 66.39541
$ws.Range("I5").Value = 0.228953097635189
$ws.Range("J5").Value = 0.228953097635189
$ws.Range("M5").Value = 0.01046566666666667
$ws.Range("N5").Value = 0.031397
$ws.Range("O5").Value = 0.007547709568116775
$ws.Range("P5").Value = 0.007547709568116775
$ws.Range("Q5").Value = 0.6948722292566666
$ws.Range("R5").Value = 6.25385006331
$ws.Range("S5").Value = 0.00172807148567109
$ws.Range("T5").Value = 0.00172807148567109
$ws.Range("G6").Value = 66.39541
$ws.Range("I6").Value = 0.228953097635189
$ws.Range("J6").Value = 0.228953097635189
$ws.Range("O6").Value = 0.8438809030711777
$ws.Range("P6").Value = 0.8438809030711776
$ws.Range("Q6").Value = 77.69103978526667
$ws.Range("R6").Value = 699.2193580674
$ws.Range("S6").Value = 0.1932091467933268
$ws.Range("T6").Value = 0.1932091467933267
$ws.Range("G7").Value = 66.39541
$ws.Range("I7").Value = 0.228953097635189
$ws.Range("J7").Value = 0.228953097635189
$ws.Range("M7").Value = 0.2060093333333333
$ws.Range("N7").Value = 0.618028
$ws.Range("O7").Value = 0.1485713873607056
$ws.Range("P7").Value = 0.1485713873607056
$ws.Range("Q7").Value = 13.67807415049333
$ws.Range("R7").Value = 123.10266735444
$ws.Range("S7").Value = 0.03401587935619111
$ws.Range("T7").Value = 0.03401587935619111
$ws.Range("G8").Value = 129.9724656666667
$ws.Range("H8").Value = 389.917397
$ws.Range("I8").Value = 0.4481875873899502
$ws.Range("J8").Value = 0.4481875873899502
$ws.Range("M8").Value = 0.01046566666666667
$ws.Range("N8").Value = 0.031397
$ws.Range("O8").Value = 0.007547709568116775
$ws.Range("P8").Value = 0.007547709568116775
$ws.Range("Q8").Value = 1.360248501512111
$ws.Range("R8").Value = 12.242236513609
$ws.Range("S8").Value = 0.003382789741654301
$ws.Range("T8").Value = 0.003382789741654301
$ws.Range("G9").Value = 129.9724656666667
$ws.Range("H9").Value = 389.917397
$ws.Range("I9").Value = 0.4481875873899502
$ws.Range("J9").Value = 0.4481875873899502
$ws.Range("O9").Value = 0.8438809030711777
$ws.Range("P9").Value = 0.8438809030711776
$ws.Range("Q9").Value = 152.0842480089844
$ws.Range("R9").Value = 1368.75823208086
$ws.Range("S9").Value = 0.3782169459919236
$ws.Range("T9").Value = 0.3782169459919235
$ws.Range("G10").Value = 129.9724656666667
$ws.Range("H10").Value = 389.917397
$ws.Range("I10").Value = 0.4481875873899502
$ws.Range("J10").Value = 0.4481875873899502
$ws.Range("M10").Value = 0.2060093333333333
$ws.Range("N10").Value = 0.618028
$ws.Range("O10").Value = 0.1485713873607056
$ws.Range("P10").Value = 0.1485713873607056
$ws.Range("Q10").Value = 26.77554100367956
$ws.Range("R10").Value = 240.979869033116
$ws.Range("S10").Value = 0.06658785165637239
$ws.Range("T10").Value = 0.06658785165637239
